# Apply updated market-price / profit figures scraped by the scheduled runner.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 20840802
$ws.Range("I62").Value = 83338090
$ws.Range("K62").Value = 83338090
$ws.Range("M62").Value = -83337466
$ws.Range("H65").Value = 20840802
$ws.Range("I65").Value = 83338090
$ws.Range("K65").Value = 416690450
$ws.Range("M65").Value = -416687330
$ws.Range("H80").Value = 2360.9656
$ws.Range("I80").Value = 658.3077
$ws.Range("J80").Value = 3744.375
$ws.Range("K80").Value = 1974.9231
$ws.Range("L80").Value = 11233.125
$ws.Range("M80").Value = -976.9231
$ws.Range("N80").Value = -13229.125
$ws.Range("H83").Value = 2360.9656
$ws.Range("I83").Value = 658.3077
$ws.Range("J83").Value = 3744.375
$ws.Range("K83").Value = 5924.7693
$ws.Range("L83").Value = 33699.375
$ws.Range("M83").Value = -932.7692999999999
$ws.Range("N83").Value = -43683.375
$ws.Range("H97").Value = 1949.8889
$ws.Range("J97").Value = 1949.8889
$ws.Range("L97").Value = 5849.6667
$ws.Range("N97").Value = -6841.6667
$ws.Range("H111").Value = 910.8461
$ws.Range("I111").Value = 859.7
$ws.Range("K111").Value = 2579.1
$ws.Range("M111").Value = 487.8999999999996
$ws.Range("H115").Value = 463.77777
$ws.Range("I115").Value = 271.75
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 815.25
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = 751.75
$ws.Range("N115").Value = -9134
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""
$ws.Range("H132").Value = 1522.2963
$ws.Range("I132").Value = 1522.2963
$ws.Range("K132").Value = 4566.8889
$ws.Range("M132").Value = -2036.8889
$ws.Range("H135").Value = 1582.1428
$ws.Range("I135").Value = 1415
$ws.Range("K135").Value = 12735
$ws.Range("M135").Value = -10200
$ws.Range("H137").Value = 14087471
$ws.Range("I137").Value = 62502520
$ws.Range("J137").Value = 3092.509
$ws.Range("K137").Value = 187507560
$ws.Range("L137").Value = 9277.527
$ws.Range("M137").Value = -187505010
$ws.Range("N137").Value = -14377.527
$ws.Range("H138").Value = 3923.8103
$ws.Range("J138").Value = 4272.04
$ws.Range("L138").Value = 12816.12
$ws.Range("N138").Value = -23096.12
$ws.Range("H141").Value = 4401
$ws.Range("I141").Value = 4401
$ws.Range("K141").Value = 13203
$ws.Range("M141").Value = -8023

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5202.4707
$ws.Range("I32").Value = 3123.3174
$ws.Range("J32").Value = 31399.8
$ws.Range("K32").Value = 3123.3174
$ws.Range("L32").Value = 31399.8
$ws.Range("M32").Value = -2836.3174
$ws.Range("N32").Value = -31973.8
$ws.Range("H61").Value = 7486.8
$ws.Range("I61").Value = 7643.1113
$ws.Range("K61").Value = 7643.1113
$ws.Range("M61").Value = -7431.1113
$ws.Range("H74").Value = 3888.25
$ws.Range("I74").Value = 2669.077
$ws.Range("K74").Value = 2669.077
$ws.Range("M74").Value = -1795.077
$ws.Range("H77").Value = 3888.25
$ws.Range("I77").Value = 2669.077
$ws.Range("K77").Value = 13345.385
$ws.Range("M77").Value = -8977.385000000002
$ws.Range("H110").Value = 5702.875
$ws.Range("I110").Value = 3801.4285
$ws.Range("K110").Value = 3801.4285
$ws.Range("M110").Value = -1756.4285
$ws.Range("H136").Value = 7486.8
$ws.Range("I136").Value = 7643.1113
$ws.Range("K136").Value = 22929.3339
$ws.Range("M136").Value = -20379.3339

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5166.081
$ws.Range("I86").Value = 4959.2905
$ws.Range("K86").Value = 4959.2905
$ws.Range("M86").Value = -3836.2905
$ws.Range("H89").Value = 5166.081
$ws.Range("I89").Value = 4959.2905
$ws.Range("K89").Value = 24796.4525
$ws.Range("M89").Value = -19180.4525
$ws.Range("H94").Value = 1129.1
$ws.Range("I94").Value = 1421.1428
$ws.Range("J94").Value = 447.66666
$ws.Range("K94").Value = 1421.1428
$ws.Range("L94").Value = 447.66666
$ws.Range("M94").Value = -970.1428000000001
$ws.Range("N94").Value = -1349.66666
$ws.Range("H107").Value = 678
$ws.Range("I107").Value = 690.5
$ws.Range("K107").Value = 690.5
$ws.Range("M107").Value = 1229.5
$ws.Range("H134").Value = 2319.814
$ws.Range("I134").Value = 1374.5
$ws.Range("K134").Value = 4123.5
$ws.Range("M134").Value = -1588.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2430.4707
$ws.Range("I16").Value = 1892.6364
$ws.Range("K16").Value = 1892.6364
$ws.Range("M16").Value = -1605.6364
$ws.Range("H58").Value = 3375.7097
$ws.Range("I58").Value = 1746.2
$ws.Range("K58").Value = 1746.2
$ws.Range("M58").Value = -1543.2
$ws.Range("H107").Value = 1450.4348
$ws.Range("I107").Value = 931.0833
$ws.Range("J107").Value = 2017
$ws.Range("K107").Value = 931.0833
$ws.Range("L107").Value = 2017
$ws.Range("M107").Value = 988.9167
$ws.Range("N107").Value = -5857
$ws.Range("H113").Value = 2430.4707
$ws.Range("I113").Value = 1892.6364
$ws.Range("K113").Value = 1892.6364
$ws.Range("M113").Value = 277.3635999999999
$ws.Range("H132").Value = 3736.3784
$ws.Range("I132").Value = 3320.0386
$ws.Range("K132").Value = 9960.1158
$ws.Range("M132").Value = -7430.1158
$ws.Range("H134").Value = 4489.6
$ws.Range("I134").Value = 3317
$ws.Range("J134").Value = 5515.625
$ws.Range("K134").Value = 9951
$ws.Range("L134").Value = 16546.875
$ws.Range("M134").Value = -7416
$ws.Range("N134").Value = -21616.875
$ws.Range("H136").Value = 3375.7097
$ws.Range("I136").Value = 1746.2
$ws.Range("K136").Value = 5238.6
$ws.Range("M136").Value = -2688.6

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2086478.5
$ws.Range("J5").Value = 4172016.5
$ws.Range("L5").Value = 12516049.5
$ws.Range("N5").Value = -12516273.5
$ws.Range("H28").Value = 720
$ws.Range("I28").Value = 720
$ws.Range("K28").Value = 2160
$ws.Range("M28").Value = -1928
$ws.Range("H134").Value = 14500470
$ws.Range("I134").Value = 7669.2856
$ws.Range("K134").Value = 23007.8568
$ws.Range("M134").Value = -17937.8568
$ws.Range("H135").Value = 2086478.5
$ws.Range("J135").Value = 4172016.5
$ws.Range("L135").Value = 37548148.5
$ws.Range("N135").Value = -37553218.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 39999.668
$ws.Range("J32").Value = 39999.668
$ws.Range("L32").Value = 39999.668
$ws.Range("N32").Value = -40591.668
$ws.Range("H126").Value = 3934.5715
$ws.Range("I126").Value = 2280.25
$ws.Range("K126").Value = 6840.75
$ws.Range("M126").Value = -4370.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 52626
$ws.Range("J112").Value = 52626
$ws.Range("L112").Value = 52626
$ws.Range("N112").Value = -55580
$ws.Range("H132").Value = 3119.7112
$ws.Range("I132").Value = 2449.3948
$ws.Range("K132").Value = 7348.1844
$ws.Range("M132").Value = -4818.1844
$ws.Range("H136").Value = 3773.36
$ws.Range("I136").Value = 2887.2307
$ws.Range("K136").Value = 8661.6921
$ws.Range("M136").Value = -6111.6921

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("H119").Value = 71723
$ws.Range("J119").Value = 71723
$ws.Range("L119").Value = 71723
$ws.Range("N119").Value = -81399
$ws.Range("H132").Value = 2638.3142
$ws.Range("I132").Value = 2134.8667
$ws.Range("K132").Value = 6404.6001
$ws.Range("M132").Value = -3874.6001
$ws.Range("H136").Value = 4543.25
$ws.Range("I136").Value = 2181
$ws.Range("K136").Value = 6543
$ws.Range("M136").Value = -3993
